$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 812.1667
$ws.Range("I34").Value = 812.1667
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 812.1667
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -609.1667
$ws.Range("N34").ClearContents()
# Row 36
$ws.Range("H36").Value = 812.1667
$ws.Range("I36").Value = 812.1667
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 812.1667
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -97.16669999999999
$ws.Range("N36").ClearContents()
# Row 40
$ws.Range("H40").Value = 1630.2858
$ws.Range("I40").Value = 1392.2
$ws.Range("J40").Value = 2225.5
$ws.Range("K40").Value = 1392.2
$ws.Range("L40").Value = 2225.5
$ws.Range("M40").Value = -1217.2
$ws.Range("N40").Value = -2575.5
# Row 132
$ws.Range("H132").Value = 27075400
$ws.Range("I132").Value = 42694004
$ws.Range("J132").Value = 3151.9333
$ws.Range("K132").Value = 128082012
$ws.Range("L132").Value = 9455.7999
$ws.Range("M132").Value = -128079482
$ws.Range("N132").Value = -14515.7999
# Row 135
$ws.Range("H135").Value = 795217.4
$ws.Range("I135").Value = 1370.8667
$ws.Range("J135").Value = 1711194.1
$ws.Range("K135").Value = 12337.8003
$ws.Range("L135").Value = 15400746.9
$ws.Range("M135").Value = -9802.800300000001
$ws.Range("N135").Value = -15405816.9
# Row 138
$ws.Range("H138").Value = 4633.717
$ws.Range("I138").Value = 6798.185
$ws.Range("J138").Value = 2386
$ws.Range("K138").Value = 20394.555
$ws.Range("L138").Value = 7158
$ws.Range("M138").Value = -15254.555
$ws.Range("N138").Value = -17438
# Row 141
$ws.Range("H141").Value = 14275
$ws.Range("I141").Value = 17208.334
$ws.Range("J141").Value = 9875
$ws.Range("K141").Value = 51625.00199999999
$ws.Range("L141").Value = 29625
$ws.Range("M141").Value = -46445.00199999999
$ws.Range("N141").Value = -39985

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7599.716
$ws.Range("I32").Value = 5441.375
$ws.Range("K32").Value = 5441.375
$ws.Range("M32").Value = -5154.375
# Row 61
$ws.Range("H61").Value = 3954.3948
$ws.Range("I61").Value = 4409.769
$ws.Range("J61").Value = 2967.75
$ws.Range("K61").Value = 4409.769
$ws.Range("L61").Value = 2967.75
$ws.Range("M61").Value = -4197.769
$ws.Range("N61").Value = -3391.75
# Row 74
$ws.Range("H74").Value = 1161.4828
$ws.Range("I74").Value = 772.4286
$ws.Range("J74").Value = 1285.2727
$ws.Range("K74").Value = 772.4286
$ws.Range("L74").Value = 1285.2727
$ws.Range("M74").Value = 101.5714
$ws.Range("N74").Value = -3033.2727
# Row 77
$ws.Range("H77").Value = 1161.4828
$ws.Range("I77").Value = 772.4286
$ws.Range("J77").Value = 1285.2727
$ws.Range("K77").Value = 3862.143
$ws.Range("L77").Value = 6426.363499999999
$ws.Range("M77").Value = 505.857
$ws.Range("N77").Value = -15162.3635
# Row 132
$ws.Range("H132").Value = 3051194.8
$ws.Range("I132").Value = 6946074.5
$ws.Range("J132").Value = 3028.0435
$ws.Range("K132").Value = 20838223.5
$ws.Range("L132").Value = 9084.130500000001
$ws.Range("M132").Value = -20835693.5
$ws.Range("N132").Value = -14144.1305
# Row 136
$ws.Range("H136").Value = 3954.3948
$ws.Range("I136").Value = 4409.769
$ws.Range("J136").Value = 2967.75
$ws.Range("K136").Value = 13229.307
$ws.Range("L136").Value = 8903.25
$ws.Range("M136").Value = -10679.307
$ws.Range("N136").Value = -14003.25

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 9023975
$ws.Range("I134").Value = 11923831
$ws.Range("J134").Value = 2199.3333
$ws.Range("K134").Value = 35771493
$ws.Range("L134").Value = 6597.999899999999
$ws.Range("M134").Value = -35768958
$ws.Range("N134").Value = -11667.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13707.056
$ws.Range("I31").Value = 1771.2858
$ws.Range("J31").Value = 21302.545
$ws.Range("K31").Value = 1771.2858
$ws.Range("L31").Value = 21302.545
$ws.Range("M31").Value = -1476.2858
$ws.Range("N31").Value = -21892.545
# Row 34
$ws.Range("H34").Value = 13707.056
$ws.Range("I34").Value = 1771.2858
$ws.Range("J34").Value = 21302.545
$ws.Range("K34").Value = 1771.2858
$ws.Range("L34").Value = 21302.545
$ws.Range("M34").Value = -1569.2858
$ws.Range("N34").Value = -21706.545
# Row 58
$ws.Range("H58").Value = 4235573.5
$ws.Range("I58").Value = 6255999
$ws.Range("J58").Value = 11048
$ws.Range("K58").Value = 6255999
$ws.Range("L58").Value = 11048
$ws.Range("M58").Value = -6255796
$ws.Range("N58").Value = -11454
# Row 132
$ws.Range("H132").Value = 6413697.5
$ws.Range("I132").Value = 11905755
$ws.Range("J132").Value = 6296.2915
$ws.Range("K132").Value = 35717265
$ws.Range("L132").Value = 18888.8745
$ws.Range("M132").Value = -35714735
$ws.Range("N132").Value = -23948.8745
# Row 134
$ws.Range("H134").Value = 8682313
$ws.Range("I134").Value = 19232508
$ws.Range("J134").Value = 2719159.2
$ws.Range("K134").Value = 57697524
$ws.Range("L134").Value = 8157477.600000001
$ws.Range("M134").Value = -57694989
$ws.Range("N134").Value = -8162547.600000001
# Row 136
$ws.Range("H136").Value = 4235573.5
$ws.Range("I136").Value = 6255999
$ws.Range("J136").Value = 11048
$ws.Range("K136").Value = 18767997
$ws.Range("L136").Value = 33144
$ws.Range("M136").Value = -18765447
$ws.Range("N136").Value = -38244

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 786.375
$ws.Range("I122").Value = 781.63635
$ws.Range("J122").Value = 796.8
$ws.Range("K122").Value = 7034.72715
$ws.Range("L122").Value = 7171.2
$ws.Range("M122").Value = -4584.72715
$ws.Range("N122").Value = -12071.2

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 181911140
$ws.Range("I11").Value = 1625
$ws.Range("J11").Value = 285859420
$ws.Range("K11").Value = 1625
$ws.Range("L11").Value = 285859420
$ws.Range("M11").Value = -1486
$ws.Range("N11").Value = -285859698
# Row 21
$ws.Range("H21").Value = 5455.4546
$ws.Range("I21").Value = 320
$ws.Range("J21").Value = 5969
$ws.Range("K21").Value = 320
$ws.Range("L21").Value = 5969
$ws.Range("M21").Value = -147
$ws.Range("N21").Value = -6315
# Row 30
$ws.Range("H30").Value = 5455.4546
$ws.Range("I30").Value = 320
$ws.Range("J30").Value = 5969
$ws.Range("K30").Value = 320
$ws.Range("L30").Value = 5969
$ws.Range("M30").Value = -215
$ws.Range("N30").Value = -6179
# Row 48
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 5000
$ws.Range("N48").Value = -5970
# Row 132
$ws.Range("H132").Value = 27031338
$ws.Range("I132").Value = 55557452
$ws.Range("J132").Value = 6598.1055
$ws.Range("K132").Value = 166672356
$ws.Range("L132").Value = 19794.3165
$ws.Range("M132").Value = -166669826
$ws.Range("N132").Value = -24854.3165

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 3010
$ws.Range("J5").Value = 3010
$ws.Range("L5").Value = 3010
$ws.Range("N5").Value = -3236
# Row 40
$ws.Range("H40").Value = 7134.8335
$ws.Range("I40").Value = 14002
$ws.Range("J40").Value = 3701.25
$ws.Range("K40").Value = 14002
$ws.Range("L40").Value = 3701.25
$ws.Range("M40").Value = -13866
$ws.Range("N40").Value = -3973.25
# Row 132
$ws.Range("H132").Value = 10529846
$ws.Range("I132").Value = 33336818
$ws.Range("J132").Value = 3552.1538
$ws.Range("K132").Value = 100010454
$ws.Range("L132").Value = 10656.4614
$ws.Range("M132").Value = -100007924
$ws.Range("N132").Value = -15716.4614

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 3907.2
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 4634
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 4634
$ws.Range("M20").Value = -760
$ws.Range("N20").Value = -5114
# Row 21
$ws.Range("H21").Value = 5555.7
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 5555.7
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 5555.7
$ws.Range("N21").Value = -6025.7
$ws.Range("M21").ClearContents()
# Row 24
$ws.Range("H24").Value = 4953
$ws.Range("I24").Value = 5000
$ws.Range("J24").Value = 4946.2856
$ws.Range("K24").Value = 5000
$ws.Range("L24").Value = 4946.2856
$ws.Range("M24").Value = -4770
$ws.Range("N24").Value = -5406.2856
# Row 26
$ws.Range("H26").Value = 4230.8887
$ws.Range("I26").Value = 1490
$ws.Range("J26").Value = 5014
$ws.Range("K26").Value = 1490
$ws.Range("L26").Value = 5014
$ws.Range("M26").Value = -1197
$ws.Range("N26").Value = -5600
# Row 28
$ws.Range("H28").Value = 5134.857
$ws.Range("J28").Value = 5134.857
$ws.Range("L28").Value = 5134.857
$ws.Range("N28").Value = -5830.857
# Row 29
$ws.Range("H29").Value = 1003146.4
$ws.Range("I29").Value = 5000000
$ws.Range("J29").Value = 3933
$ws.Range("K29").Value = 5000000
$ws.Range("L29").Value = 3933
$ws.Range("M29").Value = -4999710
$ws.Range("N29").Value = -4513
# Row 35
$ws.Range("H35").Value = 5555.7
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 5555.7
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 5555.7
$ws.Range("N35").Value = -6135.7
$ws.Range("M35").ClearContents()
# Row 38
$ws.Range("H38").Value = 8787.200000000001
$ws.Range("I38").Value = 2968
$ws.Range("J38").Value = 12666.667
$ws.Range("K38").Value = 2968
$ws.Range("L38").Value = 12666.667
$ws.Range("M38").Value = -2495
$ws.Range("N38").Value = -13612.667
# Row 40
$ws.Range("H40").Value = 4585.5713
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 4683.1665
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 4683.1665
$ws.Range("M40").Value = -3851
$ws.Range("N40").Value = -4981.1665
# Row 43
$ws.Range("H43").Value = 1686676.6
$ws.Range("J43").Value = 1686676.6
$ws.Range("L43").Value = 1686676.6
$ws.Range("N43").Value = -1686974.6
# Row 48
$ws.Range("H48").Value = 20000
$ws.Range("J48").Value = 20000
$ws.Range("L48").Value = 20000
$ws.Range("N48").Value = -21138
# Row 132
$ws.Range("H132").Value = 39721764
$ws.Range("I132").Value = 20108438
$ws.Range("J132").Value = 64238424
$ws.Range("K132").Value = 60325314
$ws.Range("L132").Value = 192715272
$ws.Range("M132").Value = -60322784
$ws.Range("N132").Value = -192720332
